$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Activate()

# Update resource names to capitalised versions, in row order (2..7)
$ws1.Range("A2").Value = "Oil"
$ws1.Range("A3").Value = "Coal"
$ws1.Range("A4").Value = "Gas"
$ws1.Range("A5").Value = "Biomass"
$ws1.Range("A6").Value = "Electricity"
$ws1.Range("A7").Value = "Hydrogen"

# Move the active selection on Feuil1 to A8 (matches recorded cursor position)
$ws1.Range("A8").Select()
